# Cryptos list data refresh — updates the Price (D) and Volume(1h) (E)
# columns for the coin rows, plus a Cronos/Fetch.AI row swap at rows
# 33-34 (coin name, link, price and volume all move together).
#
# Every Price/Volume(1h) cell in this sheet is stored as literal text
# (inline string) even when its contents look like a plain number
# (e.g. "235.62"), so a value that Excel's COM layer would otherwise
# auto-convert to a Number on assignment is written with a leading
# apostrophe first -- exactly like typing '235.62 into a cell in Excel,
# which forces text -- and the cell's style is reset to "Normal"
# afterwards so the quote-prefix formatting doesn't leave a stray style
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '96.383.93' },
    @{ Cell = 'E2'; Value = '  -0.43%  ' },
    @{ Cell = 'D3'; Value = '3.699.89' },
    @{ Cell = 'E3'; Value = '  -0.12%  ' },
    @{ Cell = 'E4'; Value = '  +0.01%  ' },
    @{ Cell = 'D5'; Value = '235.62' },
    @{ Cell = 'E5'; Value = '  -3.43%  ' },
    @{ Cell = 'E6'; Value = '  -1.98%  ' },
    @{ Cell = 'D7'; Value = '649.73' },
    @{ Cell = 'E7'; Value = '  -2.90%  ' },
    @{ Cell = 'D8'; Value = '0.427' },
    @{ Cell = 'E8'; Value = '  -0.08%  ' },
    @{ Cell = 'E9'; Value = '  +0.01%  ' },
    @{ Cell = 'E10'; Value = '  -5.41%  ' },
    @{ Cell = 'D11'; Value = '3.698.16' },
    @{ Cell = 'E11'; Value = '  -0.07%  ' },
    @{ Cell = 'D12'; Value = '0.0000308' },
    @{ Cell = 'E12'; Value = '  +16.83%  ' },
    @{ Cell = 'D13'; Value = '44.09' },
    @{ Cell = 'E13'; Value = '  -2.56%  ' },
    @{ Cell = 'D14'; Value = '0.206' },
    @{ Cell = 'E14'; Value = '  +0.34%  ' },
    @{ Cell = 'D15'; Value = '6.70' },
    @{ Cell = 'E15'; Value = '  +1.91%  ' },
    @{ Cell = 'D16'; Value = '4.388.76' },
    @{ Cell = 'E16'; Value = '  -0.10%  ' },
    @{ Cell = 'D17'; Value = '96.172.88' },
    @{ Cell = 'E17'; Value = '  -0.31%  ' },
    @{ Cell = 'D18'; Value = '8.78' },
    @{ Cell = 'E18'; Value = '  +0.19%  ' },
    @{ Cell = 'D19'; Value = '3.700.03' },
    @{ Cell = 'E19'; Value = '  -0.61%  ' },
    @{ Cell = 'E20'; Value = '  -1.04%  ' },
    @{ Cell = 'D21'; Value = '18.57' },
    @{ Cell = 'E21'; Value = '  -0.11%  ' },
    @{ Cell = 'E22'; Value = '  -7.92%  ' },
    @{ Cell = 'D23'; Value = '518.80' },
    @{ Cell = 'E23'; Value = '  +0.78%  ' },
    @{ Cell = 'D24'; Value = '3.38' },
    @{ Cell = 'E24'; Value = '  -1.55%  ' },
    @{ Cell = 'E25'; Value = '  -1.18%  ' },
    @{ Cell = 'D26'; Value = '6.89' },
    @{ Cell = 'E26'; Value = '  -0.24%  ' },
    @{ Cell = 'D27'; Value = '101.15' },
    @{ Cell = 'E27'; Value = '  -0.27%  ' },
    @{ Cell = 'D28'; Value = '13.13' },
    @{ Cell = 'E28'; Value = '  +0.68%  ' },
    @{ Cell = 'D29'; Value = '0.175' },
    @{ Cell = 'E29'; Value = '  +3.34%  ' },
    @{ Cell = 'D30'; Value = '2.99' },
    @{ Cell = 'E30'; Value = '  -2.70%  ' },
    @{ Cell = 'D31'; Value = '12.06' },
    @{ Cell = 'E31'; Value = '  +0.12%  ' },
    @{ Cell = 'E32'; Value = '  +0.21%  ' },
    @{ Cell = 'B33'; Value = 'Fetch.AI' },
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet' },
    @{ Cell = 'D33'; Value = '1.86' },
    @{ Cell = 'E33'; Value = '  +6.56%  ' },
    @{ Cell = 'B34'; Value = 'Cronos' },
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' },
    @{ Cell = 'D34'; Value = '0.185' },
    @{ Cell = 'E34'; Value = '  -0.38%  ' },
    @{ Cell = 'D35'; Value = '0.998' },
    @{ Cell = 'E35'; Value = '  +2.17%  ' },
    @{ Cell = 'D36'; Value = '32.15' },
    @{ Cell = 'E36'; Value = '  -4.00%  ' },
    @{ Cell = 'D37'; Value = '646.87' },
    @{ Cell = 'E37'; Value = '  +5.09%  ' },
    @{ Cell = 'D38'; Value = '0.584' },
    @{ Cell = 'E38'; Value = '  -1.65%  ' },
    @{ Cell = 'D39'; Value = '8.79' },
    @{ Cell = 'E39'; Value = '  +0.30%  ' },
    @{ Cell = 'E40'; Value = '  +0.04%  ' },
    @{ Cell = 'D41'; Value = '6.82' },
    @{ Cell = 'E41'; Value = '  +11.34%  ' },
    @{ Cell = 'E42'; Value = '  +4.48%  ' },
    @{ Cell = 'D43'; Value = '40.70' },
    @{ Cell = 'E43'; Value = '  -4.67%  ' },
    @{ Cell = 'E44'; Value = '  -0.38%  ' },
    @{ Cell = 'D45'; Value = '0.955' },
    @{ Cell = 'E45'; Value = '  -1.53%  ' },
    @{ Cell = 'D46'; Value = '0.0448' },
    @{ Cell = 'E46'; Value = '  +0.95%  ' },
    @{ Cell = 'E47'; Value = '  +1.77%  ' },
    @{ Cell = 'E48'; Value = '  -0.06%  ' },
    @{ Cell = 'E49'; Value = '  -1.54%  ' },
    @{ Cell = 'D50'; Value = '8.44' },
    @{ Cell = 'E50'; Value = '  -1.99%  ' },
    @{ Cell = 'E51'; Value = '  +2.26%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $value = $u.Value

    if ($value -match '^-?\d+(\.\d+)?$') {
        # Plain-number-looking text (e.g. a Price cell) -- force text so
        # it isn't coerced into a Number, then strip the leftover
        # quote-prefix style.
        $range.Value = "'" + $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
